$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Userstory")

# --- Fill in "Datum" (E) and "Status" (F) for the existing user stories
#     (rows 2-9), and add the two new user stories (rows 10-11). ---

$ws.Range("E2").Value = 42079
$ws.Range("F2").Value = "In Progress"

$ws.Range("E3").Value = 42079
$ws.Range("F3").Value = "In Progress"

$ws.Range("E4").Value = 42079
$ws.Range("F4").Value = "Finished"

$ws.Range("E5").Value = 42079
$ws.Range("F5").Value = "In Progress"

$ws.Range("E6").Value = 42079
$ws.Range("F6").Value = "In Progress"

$ws.Range("E7").Value = 42079
$ws.Range("F7").Value = "In Progress"

$ws.Range("E8").Value = 42079
$ws.Range("F8").Value = "In Progress"

$ws.Range("E9").Value = 42079
$ws.Range("F9").Value = "In Progress"

$ws.Range("B10").Value = "Webclient coderen."
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "Matthias"
$ws.Range("E10").Value = 42080
$ws.Range("F10").Value = "In Progress"

$ws.Range("B11").Value = "Box fixing en fixed nav + footer"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Matthias"
$ws.Range("E11").Value = 42080
$ws.Range("F11").Value = "Finished"

# E2 already carries the short-date style; copy it onto the rest of the
# "Datum" column so every new date cell renders (and is stored) the same way.
$ws.Range("E2").Copy()
$ws.Range("E3:E11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Match the selection left behind by the author.
$ws.Range("G11").Select()
